# "feat: basic app working / moved database folder / generates history of menus"
#
# The "Pratos" sheet (sheet1) had a week's worth of planned-meal rows
# (rows 2-10). The app now generates/consumes that history itself, so the
# previously hand-entered rows are cleared out, leaving only the header row
# and the two rows whose "Carne" cell (column B) keeps its special font
# style (s="1") as a template/placeholder for new entries.
#
# Clearing the cells (rather than deleting rows/using the ListObject API)
# keeps row numbers, the table's declared range, and the data-validation
# ranges exactly as they were - only the cell contents disappear.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pratos")

# Row 2: clear all three cells. ClearContents() only removes the value
# (and the t="s" type attribute) - the cell's style (s="1" on B2) stays.
$ws.Range("A2:C2").ClearContents()

# Rows 3-5: fully cleared (no styled cells to preserve there).
$ws.Range("A3:C5").ClearContents()

# Row 6: same treatment as row 2 - B6 keeps its style after clearing.
$ws.Range("A6:C6").ClearContents()

# Rows 7-10: fully cleared.
$ws.Range("A7:C10").ClearContents()

# Move the active selection back onto the sheet (it was parked on E12,
# well outside the shrunk used range).
$ws.Range("A2").Select() | Out-Null
